$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.438.03'
$ws.Range("D3").Value = '2.069.81'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.09'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("E6").Value = '  +1.37%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.19'
$ws.Range("E8").Value = '  -1.83%  '
$ws.Range("E9").Value = '  +2.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0774'
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").Value = '2.373.11'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.40'
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.89'
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.775'
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '2.069.66'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '37.354.29'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.18'
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.58'
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.76'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.39'
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("E28").Value = '  -7.01%  '
$ws.Range("E29").Value = '  +1.06%  '
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("E31").Value = '  -1.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.52'
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0615'
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("E35").Value = '  -2.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.38'
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.60'
$ws.Range("E39").Value = '  -4.22%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").Value = '  -3.35%  '
$ws.Range("D42").Value = '1.484.27'
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.07'
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("E46").Value = '  -4.88%  '
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.17'
$ws.Range("E48").Value = '  -4.63%  '
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.94'
$ws.Range("E51").Value = '  +4.76%  '
